$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I0 (col I) and IF (col J)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style from H1 onto the two new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data for columns I (I0) and J (IF), for rows 2..69
$iVals = @(6,6,6,6,7,7,7,6,6,6,5,6,3,5,7,7,4,7,8,6,8,8,6,4,6,6,4,8,5,5,5,7,6,5,9,5,7,8,6,9,6,7,7,8,8,7,8,6,9,8,6,8,4,4,8,6,6,4,9,6,9,6,9,7,5,2,8,3)
$jVals = @(7,6,7,7,8,7,8,6,7,7,6,6,4,5,7,7,4,7,8,6,8,8,7,5,7,6,5,8,6,6,6,7,6,5,9,5,7,8,6,9,6,7,8,8,8,8,8,6,9,8,7,9,4,6,9,7,7,5,9,7,9,6,9,7,5,3,8,3)

for ($idx = 0; $idx -lt $iVals.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$idx]
    $ws.Cells.Item($row, 10).Value = $jVals[$idx]
}
